# ADD results from server
# Updates row 2 (data row) values on each year sheet (2025, 2030, 2035, 2040, 2045, 2050)
# to reflect the latest server results.

$wb = $excel.ActiveWorkbook

# Sheet 1: "2025"
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 277.8902526399997
$ws.Range("E2").Value = 29092.72506141524
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 14865.25872276
$ws.Range("L2").Value = 50912.59821312752
$ws.Range("M2").Value = 11247.09127927
$ws.Range("N2").Value = 7270.39941619098
$ws.Range("O2").Value = 6890.515200515631

# Sheet 2: "2030"
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5506.32916570769
$ws.Range("E2").Value = 56005.10427174018
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 36859.48338500719
$ws.Range("L2").Value = 73837.13070482886
$ws.Range("M2").Value = 21806.186745563
$ws.Range("N2").Value = 10959.27269082414
$ws.Range("O2").Value = 9417.058399985341

# Sheet 3: "2035"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2232.402753895485
$ws.Range("B2").Value = 6951.177638494421
$ws.Range("E2").Value = 67289.88340938435
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 53206.8045886644
$ws.Range("L2").Value = 73837.13070482886
$ws.Range("M2").Value = 27593.46631997451
$ws.Range("N2").Value = 15985.96622368204
$ws.Range("O2").Value = 15294.59034761512

# Sheet 4: "2040"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2232.402753895485
$ws.Range("B2").Value = 6951.177638494421
$ws.Range("E2").Value = 67289.88340938435
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 53206.8045886644
$ws.Range("L2").Value = 73837.13070482886
$ws.Range("M2").Value = 27593.46631997451
$ws.Range("N2").Value = 15985.96622368204
$ws.Range("O2").Value = 15294.59034761512

# Sheet 5: "2045"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 2232.402753895485
$ws.Range("B2").Value = 6951.177638494421
$ws.Range("E2").Value = 67289.88340938435
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 53206.8045886644
$ws.Range("L2").Value = 73837.13070482886
$ws.Range("M2").Value = 27593.46631997451
$ws.Range("N2").Value = 15985.96622368204
$ws.Range("O2").Value = 15294.59034761512

# Sheet 6: "2050"
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 2232.402753895485
$ws.Range("B2").Value = 6951.177638494421
$ws.Range("E2").Value = 67289.88340938435
$ws.Range("G2").Value = 8095.925712662093
$ws.Range("I2").Value = 53206.8045886644
$ws.Range("L2").Value = 73837.13070482886
$ws.Range("M2").Value = 27593.46631997451
$ws.Range("N2").Value = 15985.96622368204
$ws.Range("O2").Value = 15294.59034761512
